# boosted regression tree - re-organize & new for all plant species richness
#
# Adds a new "All plant richness" BRT-output block (rows 167-201), mirroring
# the structure of the first block on the sheet (rows 1-40), extends the
# print area / dimension, and adds a page break before the new block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1) Copy cell formatting for the new block from the analogous cells of
#    the existing first block (rows 1-40), which has an identical layout
#    (title / header / importance table / deviance table / interactions).
# -----------------------------------------------------------------------
$ws.Range("A1:C1").Copy()
$ws.Range("A167:C167").PasteSpecial(-4122)

$ws.Range("A3:B3").Copy()
$ws.Range("A169:B169").PasteSpecial(-4122)

$ws.Range("A4:B23").Copy()
$ws.Range("A170:B189").PasteSpecial(-4122)

$ws.Range("A25:B25").Copy()
$ws.Range("A190:B190").PasteSpecial(-4122)

$ws.Range("A26:B27").Copy()
$ws.Range("A191:B192").PasteSpecial(-4122)

$ws.Range("A28:B28").Copy()
$ws.Range("A193:B193").PasteSpecial(-4122)

$ws.Range("A29:C29").Copy()
$ws.Range("A194:C194").PasteSpecial(-4122)

$ws.Range("A30").Copy()
$ws.Range("A195").PasteSpecial(-4122)

$ws.Range("A31:C31").Copy()
$ws.Range("A196:C196").PasteSpecial(-4122)

$ws.Range("A32:C35").Copy()
$ws.Range("A197:C200").PasteSpecial(-4122)

$ws.Range("A40:C40").Copy()
$ws.Range("A201:C201").PasteSpecial(-4122)

$ws.Range("A1").Select()

# -----------------------------------------------------------------------
# 2) Title
# -----------------------------------------------------------------------
$ws.Range("A167").Value = "All plant richness"

# -----------------------------------------------------------------------
# 3) Variable importance table
# -----------------------------------------------------------------------
$ws.Range("A169").Value = "Variable"
$ws.Range("B169").Value = "Importance"

$ws.Range("A170").Value = "regional_watershed"
$ws.Range("B170").Value = 50.141071009999997

$ws.Range("A171").Value = "surfacearea_ha"
$ws.Range("B171").Value = 19.668259450000001

$ws.Range("A172").Value = "boatlaunch"
$ws.Range("B172").Value = 6.7384688600000002

$ws.Range("A173").Value = "secchi_avg"
$ws.Range("B173").Value = 6.1890596499999999

$ws.Range("A174").Value = "shoreline_development"
$ws.Range("B174").Value = 5.6525875799999996

$ws.Range("A175").Value = "COND_avg"
$ws.Range("B175").Value = 2.8722203999999998

$ws.Range("A176").Value = "ALK_avg"
$ws.Range("B176").Value = 2.31918523

$ws.Range("A177").Value = "nearest_LM"
$ws.Range("B177").Value = 1.5034017200000001

$ws.Range("A178").Value = "depth_max_m"
$ws.Range("B178").Value = 1.2325731499999999

$ws.Range("A179").Value = "nearest_W"
$ws.Range("B179").Value = 0.90205245000000001

$ws.Range("A180").Value = "dist_waterfowl"
$ws.Range("B180").Value = 0.65937913999999997

$ws.Range("A181").Value = "longitude"
$ws.Range("B181").Value = 0.47221640999999998

$ws.Range("A182").Value = "TOTP_avg"
$ws.Range("B182").Value = 0.38907717000000003

$ws.Range("A183").Value = "waterbodies_5km"
$ws.Range("B183").Value = 0.363954

$ws.Range("A184").Value = "PH_avg"
$ws.Range("B184").Value = 0.28998081999999997

$ws.Range("A185").Value = "nearest_SP"
$ws.Range("B185").Value = 0.26705433000000001

$ws.Range("A186").Value = "waterbodies_10km"
$ws.Range("B186").Value = 0.22551354000000001

$ws.Range("A187").Value = "latitude"
$ws.Range("B187").Value = 0.09360299

$ws.Range("A188").Value = "major_watershed"
$ws.Range("B188").Value = 0.01297081

$ws.Range("A189").Value = "waterbodies_1km"
$ws.Range("B189").Value = 0.00737127

# -----------------------------------------------------------------------
# 4) Null / residual deviance + % deviance explained
# -----------------------------------------------------------------------
$ws.Range("A191").Value = "Null deviance"
$ws.Range("B191").Value = 863.52110000000005

$ws.Range("A192").Value = "Resid deviance"
$ws.Range("B192").Value = 374.85719999999998

$ws.Range("A194").Value = "% Deviance Explained"
$ws.Range("B194").Value = "(Dev.null - Dev.resid) / Dev.null * 100"
$ws.Range("C194").Formula = "=(B191-B192)/B191*100"

# -----------------------------------------------------------------------
# 5) Interaction table
# -----------------------------------------------------------------------
$ws.Range("A196").Value = "Variable1"
$ws.Range("B196").Value = "Variable2"
$ws.Range("C196").Value = "Interaction.size"

$ws.Range("A197").Value = "regional_watershed"
$ws.Range("B197").Value = "surfacearea_ha"
$ws.Range("C197").Value = 0.14000000000000001

$ws.Range("A198").Value = "boatlaunch"
$ws.Range("B198").Value = "regional_watershed"
$ws.Range("C198").Value = 0.13

$ws.Range("A199").Value = "regional_watershed"
$ws.Range("B199").Value = "secchi_avg"
$ws.Range("C199").Value = 0.03

$ws.Range("A200").Value = "regional_watershed"
$ws.Range("B200").Value = "shoreline_development"
$ws.Range("C200").Value = 0.02

$ws.Range("A201").Value = "Not reporting 0 interactions"

# -----------------------------------------------------------------------
# 6) Print area, page break, view/selection
# -----------------------------------------------------------------------
$ws.PageSetup.PrintArea = "`$A`$1:`$D`$201"

$ws.HPageBreaks.Add($ws.Range("A167"))

$excel.ActiveWindow.Zoom = 60
$ws.Range("C193").Select()
